# Apply the "Compound List" update:
#  - rename/reorder the density/Vm/MW columns
#  - insert a new "benzene" row after ethanol
#  - swap the Vm / MW values for the pre-existing compounds

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3 (pushes the existing rows 3-5 down to 4-6)
$ws.Rows.Item(3).Insert()

# --- Header row ---
$ws.Range("C1").Value = "ρ (g/ml)"
$ws.Range("D1").Value = "Vm (ml/mol)"
$ws.Range("E1").Value = "MW (g/mol)"

# --- Row 2: ethanol -- swap Vm/MW (D/E) ---
$ws.Range("D2").Value = 57.93
$ws.Range("E2").Value = 46.042

# --- Row 3 (new): benzene ---
$ws.Range("A3").Value = "benzene"
$ws.Range("B3").Value = "c1ccccc1"
$ws.Range("C3").Value = 0.886
$ws.Range("D3").Value = 88.08999999999999
$ws.Range("E3").Value = 78.047

# --- Row 4: 2-Isobutyl-4-methyl-1,3-dioxolane -- swap Vm/MW (D/E) ---
$ws.Range("D4").Value = 157.17
$ws.Range("E4").Value = 144.115

# --- Row 5: 2-Ethyl-4-methyl-1,3-dioxolane -- swap Vm/MW (D/E) ---
$ws.Range("D5").Value = 123.13
$ws.Range("E5").Value = 116.084

# --- Row 6: 3,3-Dimethyloxetane -- swap Vm/MW (D/E) ---
$ws.Range("D6").Value = 101.04
$ws.Range("E6").Value = 86.07299999999999
